{"js": "// \"Joel 8====D, Justin weekend of 11/17/17\" -> \"Joel, Justin weekend of 11/17/17\"\n// (drop the \" 8====D\" aside), and relocate the stray _GoBack bookmark from the\n// end of the following \"Start with fastq files...\" paragraph to sit right\n// after \"Joel\" (where the cursor actually was when the text got edited).\n\nconst body = context.document.body;\n\n// 1. Remove \" 8====D\" from the \"Joel 8====D, Justin weekend of 11/17/17\" line.\nconst toDelete = body.search(\" 8====D\", { matchCase: true });\ntoDelete.load(\"items\");\nawait context.sync();\n\nif (toDelete.items.length > 0) {\n  toDelete.items[0].insertText(\"\", \"Replace\");\n}\nawait context.sync();\n\n// 2. Move the _GoBack bookmark: delete it from wherever it currently sits\n//    (end of the \"Start with fastq files...\" paragraph) and re-add it right\n//    after \"Joel\" in the paragraph above.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst joelHits = body.search(\"Joel\", { matchCase: true });\njoelHits.load(\"items\");\nawait context.sync();\n\nif (joelHits.items.length > 0) {\n  const afterJoel = joelHits.items[0].getRange(\"End\");\n  afterJoel.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# \"Joel 8====D, Justin weekend of 11/17/17\" -> \"Joel, Justin weekend of 11/17/17\"\n# (drop the \" 8====D\" aside), and relocate the stray _GoBack bookmark from the\n# end of the following \"Start with fastq files...\" paragraph to sit right\n# after \"Joel\" (where the cursor actually was when the text got edited).\n\n$d = $word.ActiveDocument\n\n# 1. Remove \" 8====D\" from the \"Joel 8====D, Justin weekend of 11/17/17\" line.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\" 8====D\")\nif ($found) {\n    $rng.Text = \"\"\n}\n\n# 2. Move the _GoBack bookmark: delete it from wherever it currently sits\n#    (end of the \"Start with fastq files...\" paragraph) and re-add it right\n#    after \"Joel\" in the paragraph above.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$joel = $d.Content\n$joel.Find.ClearFormatting()\n$joel.Find.Execute(\"Joel\") | Out-Null\n$afterJoel = $d.Range($joel.End, $joel.End)\n$d.Bookmarks.Add(\"_GoBack\", $afterJoel)\n"}
